# Apply updated Betfair Back/Lay odds for 2025-11-17
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = 2.92   # F2
$ws.Cells.Item(2, 7).Value = 4.8   # G2
$ws.Cells.Item(2, 8).Value = 1.92   # H2
$ws.Cells.Item(2, 9).Value = 2.68   # I2
$ws.Cells.Item(2, 10).Value = 2.88   # J2
$ws.Cells.Item(2, 11).Value = 7.4   # K2
$ws.Cells.Item(2, 12).Value = 1.01   # L2
$ws.Cells.Item(2, 13).Value = 1.01   # M2
$ws.Cells.Item(2, 14).Value = 1.58   # N2
$ws.Cells.Item(2, 15).Value = 1.01   # O2
$ws.Cells.Item(2, 16).Value = 1.57   # P2
$ws.Cells.Item(2, 17).Value = 1.93   # Q2
$ws.Cells.Item(2, 18).Value = 1.18   # R2
$ws.Cells.Item(2, 19).Value = 1.94   # S2
$ws.Cells.Item(2, 20).Value = 1.01   # T2
$ws.Cells.Item(2, 21).Value = 1.01   # U2
$ws.Cells.Item(2, 22).Value = 1.59   # V2
$ws.Cells.Item(2, 23).Value = 1.26   # W2
$ws.Cells.Item(2, 24).Value = 1000   # X2
$ws.Cells.Item(2, 25).Value = 1000   # Y2
$ws.Cells.Item(2, 26).Value = 1000   # Z2
$ws.Cells.Item(2, 27).Value = 1000   # AA2
$ws.Cells.Item(2, 28).Value = 1000   # AB2
$ws.Cells.Item(2, 29).Value = 1000   # AC2
$ws.Cells.Item(2, 30).Value = 1000   # AD2
$ws.Cells.Item(2, 31).Value = 1000   # AE2
$ws.Cells.Item(2, 32).Value = 1000   # AF2
$ws.Cells.Item(2, 33).Value = 1000   # AG2
$ws.Cells.Item(2, 34).Value = 1000   # AH2
$ws.Cells.Item(2, 35).Value = 1000   # AI2
$ws.Cells.Item(2, 36).Value = 1000   # AJ2
$ws.Cells.Item(2, 37).Value = 1000   # AK2
$ws.Cells.Item(2, 38).Value = 1000   # AL2
$ws.Cells.Item(2, 39).Value = 1000   # AM2
$ws.Cells.Item(2, 40).Value = 1000   # AN2
$ws.Cells.Item(2, 41).Value = 1000   # AO2
$ws.Cells.Item(3, 6).Value = 2.32   # F3
$ws.Cells.Item(3, 7).Value = 3.3   # G3
$ws.Cells.Item(3, 8).Value = 2.66   # H3
$ws.Cells.Item(3, 9).Value = 4   # I3
$ws.Cells.Item(3, 10).Value = 2.62   # J3
$ws.Cells.Item(3, 11).Value = 5.2   # K3
$ws.Cells.Item(3, 12).Value = 1.01   # L3
$ws.Cells.Item(3, 13).Value = 1.01   # M3
$ws.Cells.Item(3, 14).Value = 1.35   # N3
$ws.Cells.Item(3, 15).Value = 1.01   # O3
$ws.Cells.Item(3, 16).Value = 1.35   # P3
$ws.Cells.Item(3, 17).Value = 2.52   # Q3
$ws.Cells.Item(3, 18).Value = 1.18   # R3
$ws.Cells.Item(3, 19).Value = 2.52   # S3
$ws.Cells.Item(3, 20).Value = 1.01   # T3
$ws.Cells.Item(3, 21).Value = 1.01   # U3
$ws.Cells.Item(3, 22).Value = 1.33   # V3
$ws.Cells.Item(3, 23).Value = 1.43   # W3
$ws.Cells.Item(3, 24).Value = 1000   # X3
$ws.Cells.Item(3, 25).Value = 1000   # Y3
$ws.Cells.Item(3, 26).Value = 1000   # Z3
$ws.Cells.Item(3, 27).Value = 1000   # AA3
$ws.Cells.Item(3, 28).Value = 1000   # AB3
$ws.Cells.Item(3, 29).Value = 1000   # AC3
$ws.Cells.Item(3, 30).Value = 1000   # AD3
$ws.Cells.Item(3, 31).Value = 1000   # AE3
$ws.Cells.Item(3, 32).Value = 1000   # AF3
$ws.Cells.Item(3, 33).Value = 1000   # AG3
$ws.Cells.Item(3, 34).Value = 1000   # AH3
$ws.Cells.Item(3, 35).Value = 1000   # AI3
$ws.Cells.Item(3, 36).Value = 1000   # AJ3
$ws.Cells.Item(3, 37).Value = 1000   # AK3
$ws.Cells.Item(3, 38).Value = 1000   # AL3
$ws.Cells.Item(3, 39).Value = 1000   # AM3
$ws.Cells.Item(3, 40).Value = 1000   # AN3
$ws.Cells.Item(3, 41).Value = 1000   # AO3
$ws.Cells.Item(4, 6).Value = 3.3   # F4
$ws.Cells.Item(4, 7).Value = 3.45   # G4
$ws.Cells.Item(4, 8).Value = 2.74   # H4
$ws.Cells.Item(4, 9).Value = 2.9   # I4
$ws.Cells.Item(4, 11).Value = 2.96   # K4
$ws.Cells.Item(4, 14).Value = 1.94   # N4
$ws.Cells.Item(4, 19).Value = 9.8   # S4
$ws.Cells.Item(4, 21).Value = 1.47   # U4
$ws.Cells.Item(4, 22).Value = 1.52   # V4
$ws.Cells.Item(4, 23).Value = 1.41   # W4
$ws.Cells.Item(4, 36).Value = 980   # AJ4
$ws.Cells.Item(5, 8).Value = 3.85   # H5
$ws.Cells.Item(5, 11).Value = 3.4   # K5
$ws.Cells.Item(5, 25).Value = 13.5   # Y5
$ws.Cells.Item(6, 9).Value = 5.2   # I6
$ws.Cells.Item(6, 10).Value = 3.1   # J6
$ws.Cells.Item(6, 11).Value = 3.35   # K6
$ws.Cells.Item(6, 12).Value = 1.46   # L6
$ws.Cells.Item(6, 20).Value = 1.92   # T6
$ws.Cells.Item(6, 21).Value = 1.78   # U6
$ws.Cells.Item(6, 24).Value = 1000   # X6
$ws.Cells.Item(6, 25).Value = 980   # Y6
$ws.Cells.Item(6, 27).Value = 1000   # AA6
$ws.Cells.Item(6, 28).Value = 980   # AB6
$ws.Cells.Item(6, 29).Value = 980   # AC6
$ws.Cells.Item(6, 31).Value = 1000   # AE6
$ws.Cells.Item(6, 35).Value = 130   # AI6
$ws.Cells.Item(6, 38).Value = 1000   # AL6
$ws.Cells.Item(6, 40).Value = 1000   # AN6
$ws.Cells.Item(7, 6).Value = 2.26   # F7
$ws.Cells.Item(7, 7).Value = 2.46   # G7
$ws.Cells.Item(7, 8).Value = 3.7   # H7
$ws.Cells.Item(7, 9).Value = 4.2   # I7
$ws.Cells.Item(7, 10).Value = 3.05   # J7
$ws.Cells.Item(7, 20).Value = 2   # T7
$ws.Cells.Item(7, 21).Value = 1.65   # U7
$ws.Cells.Item(7, 22).Value = 1.32   # V7
$ws.Cells.Item(7, 25).Value = 9.8   # Y7
$ws.Cells.Item(7, 26).Value = 980   # Z7
$ws.Cells.Item(7, 28).Value = 8   # AB7
